# Services.xlsx update
# - "contact us"/"payment " columns collapse into a single "actionType" column:
#     H1: "contact us"  -> "actionType"
#     H2/H3 (HVAC rows): "Yes" -> "both"
#     H4/H5/H6 (fiber + install rows): "Yes" -> "contact"
#     column I ("Yes"/"No"/"payment ") is removed entirely
# - the "price" / "sales price" values were transposed for the two fiber rows
#     row 4 (fiber-sedan): C<->D swap (999/800 -> 800/999)
#     row 6 (fiber-suv):   C<->D swap (1200/1000 -> 1000/1200)
# - selection moves to I1 (now the first empty column after the delete)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the old "contact us" / "payment" pair of columns into one
# "actionType" column (H), recording whether the action is "both" (call or
# book online) or just "contact".
$ws.Range("H1").Value = "actionType"
$ws.Range("H2").Value = "both"
$ws.Range("H3").Value = "both"
$ws.Range("H4").Value = "contact"
$ws.Range("H5").Value = "contact"
$ws.Range("H6").Value = "contact"

# Swap the price/sales-price figures that were entered in the wrong columns
# for the two fiber-optic rows.
$c4 = $ws.Range("C4").Value()
$d4 = $ws.Range("D4").Value()
$ws.Range("C4").Value = $d4
$ws.Range("D4").Value = $c4

$c6 = $ws.Range("C6").Value()
$d6 = $ws.Range("D6").Value()
$ws.Range("C6").Value = $d6
$ws.Range("D6").Value = $c6

# The old "payment" column (I) is no longer needed now that it has been
# folded into "actionType".
$ws.Columns("I").Delete()

# Leave the selection where the workbook was last saved.
$ws.Range("I1").Select()
